# 7.8 History Card & Advanced Story
#
# Rework the opening lines of the Study scene on Sheet1:
#   - B11 (the old "Click on any area you find suspicious..." hint) becomes a
#     green-colored, parenthesized hint line, and its row grows to a 2-line height.
#   - B5 ("Please, both of you, come in.") is shortened to "Please come in."
#   - B4 (the "Steward He skillfully unlocked..." line) becomes a green-colored,
#     parenthesized narration line describing the Butler unlocking the study,
#     and its row grows to a 2-line height.
#   - The sheet's saved selection moves from B20 to B4.
#
# NOTE: the order in which new strings are first written matters because it
# determines the order they are appended to the shared-string table, so we
# touch B11, then B5, then B4 to reproduce the expected table layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hint shown just before investigation begins.
$ws.Range("B11").Value = " <color=#00CC00>(Click on any area you find suspicious to gather clues.)</color>"
$ws.Rows.Item(11).RowHeight = 34

# Shorten the butler's invitation line.
$ws.Range("B5").Value = "Please come in."

# Narration describing the butler unlocking the study door.
$ws.Range("B4").Value = " <color=#00CC00>(Butler He skillfully unlocked the door to the study.)</color>"
$ws.Rows.Item(4).RowHeight = 34

# Update the sheet's saved selection/active cell.
$ws.Range("B4").Select()
